$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.365.42"
$ws.Range("E2").Value = "  +12.76%  "
$ws.Range("D3").Value = "1.822.76"
$ws.Range("E3").Value = "  +7.97%  "
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").Value = "'232.11"
$ws.Range("E5").Value = "  +4.99%  "
$ws.Range("D6").Value = "'0.548"
$ws.Range("E6").Value = "  +4.78%  "
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("D8").Value = "'31.57"
$ws.Range("E8").Value = "  +3.67%  "
$ws.Range("D9").Value = "'46.16"
$ws.Range("E9").Value = "  +1.99%  "
$ws.Range("D10").Value = "'0.284"
$ws.Range("E10").Value = "  +7.12%  "
$ws.Range("D11").Value = "'0.0680"
$ws.Range("E11").Value = "  +9.10%  "
$ws.Range("E12").Value = "  +3.30%  "
$ws.Range("D13").Value = "2.085.42"
$ws.Range("E13").Value = "  +7.97%  "
$ws.Range("D14").Value = "1.817.27"
$ws.Range("E14").Value = "  +7.70%  "
$ws.Range("D15").Value = "'0.646"
$ws.Range("E15").Value = "  +4.77%  "
$ws.Range("D16").Value = "34.378.96"
$ws.Range("E16").Value = "  +12.48%  "
$ws.Range("D17").Value = "'10.32"
$ws.Range("E17").Value = "  -3.99%  "
$ws.Range("E18").Value = "  +8.28%  "
$ws.Range("D19").Value = "'70.82"
$ws.Range("E19").Value = "  +7.21%  "
$ws.Range("D20").Value = "'260.77"
$ws.Range("E20").Value = "  +5.55%  "
$ws.Range("D21").Value = "0.0₃0752"
$ws.Range("E21").Value = "  +4.43%  "
$ws.Range("E22").Value = "  +0.34%  "
$ws.Range("D23").Value = "'10.51"
$ws.Range("E23").Value = "  +3.10%  "
$ws.Range("E24").Value = "  +2.75%  "
$ws.Range("D25").Value = "'2.19"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("D26").Value = "'161.28"
$ws.Range("E26").Value = "  +1.77%  "
$ws.Range("D27").Value = "'16.85"
$ws.Range("E27").Value = "  +5.80%  "
$ws.Range("E28").Value = "  +4.81%  "
$ws.Range("D29").Value = "'7.15"
$ws.Range("E29").Value = "  +5.51%  "
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").Value = "'3.84"
$ws.Range("E31").Value = "  +9.46%  "
$ws.Range("D32").Value = "'1.22"
$ws.Range("E32").Value = "  +7.32%  "
$ws.Range("D33").Value = "'0.0517"
$ws.Range("E33").Value = "  +2.97%  "
$ws.Range("E34").Value = "  +7.98%  "
$ws.Range("D35").Value = "1.589.75"
$ws.Range("E35").Value = "  +5.33%  "
$ws.Range("E36").Value = "  +5.95%  "
$ws.Range("E37").Value = "  +2.68%  "
$ws.Range("D38").Value = "'85.67"
$ws.Range("E38").Value = "  +8.15%  "
$ws.Range("E39").Value = "  +5.31%  "
$ws.Range("D40").Value = "'0.630"
$ws.Range("E40").Value = "  +7.64%  "
$ws.Range("D41").Value = "'2.78"
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("D42").Value = "'2.38"
$ws.Range("E42").Value = "  +2.64%  "
$ws.Range("D43").Value = "'0.922"
$ws.Range("E43").Value = "  +8.23%  "
$ws.Range("E44").Value = "  +6.07%  "
$ws.Range("E45").Value = "  +3.73%  "
$ws.Range("E46").Value = "  +6.80%  "
$ws.Range("D47").Value = "1.976.53"
$ws.Range("E47").Value = "  +8.16%  "
$ws.Range("D48").Value = "'53.59"
$ws.Range("E48").Value = "  +2.57%  "
$ws.Range("D49").Value = "'5.73"
$ws.Range("E49").Value = "  +5.69%  "
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'11.36"
$ws.Range("E51").Value = "  +21.85%  "
